# Update the cryptos price/volume snapshot (column D = Price, column E = Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.332.27'
$ws.Range('E2').Value = '  -3.07%  '
$ws.Range('D3').Value = '2.454.97'
$ws.Range('E3').Value = '  -3.25%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''310.86'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '''93.62'
$ws.Range('E6').Value = '  -6.96%  '
$ws.Range('E7').Value = '  -3.55%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.505'
$ws.Range('E9').Value = '  -4.52%  '
$ws.Range('D10').Value = '''33.21'
$ws.Range('E10').Value = '  -7.51%  '
$ws.Range('D11').Value = '''0.0780'
$ws.Range('E11').Value = '  -3.10%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '''6.91'
$ws.Range('E13').Value = '  -6.01%  '
$ws.Range('D14').Value = '2.835.30'
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('D15').Value = '2.439.78'
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').Value = '''14.35'
$ws.Range('E16').Value = '  -9.71%  '
$ws.Range('D17').Value = '''0.786'
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('D18').Value = '41.340.95'
$ws.Range('E18').Value = '  -3.08%  '
$ws.Range('D19').Value = '''6.32'
$ws.Range('E19').Value = '  -6.52%  '
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  -4.44%  '
$ws.Range('D21').Value = '''11.46'
$ws.Range('E21').Value = '  -6.35%  '
$ws.Range('D22').Value = '''68.43'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').Value = '''238.08'
$ws.Range('E23').Value = '  -2.34%  '
$ws.Range('D24').Value = '''2.76'
$ws.Range('E24').Value = '  -4.76%  '
$ws.Range('D25').Value = '''1.92'
$ws.Range('E25').Value = '  -6.17%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '''24.64'
$ws.Range('E27').Value = '  -5.38%  '
$ws.Range('E28').Value = '  -5.97%  '
$ws.Range('D29').Value = '''9.67'
$ws.Range('E29').Value = '  -4.83%  '
$ws.Range('D30').Value = '''36.22'
$ws.Range('E30').Value = '  -7.82%  '
$ws.Range('D31').Value = '''152.56'
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').Value = '''5.58'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -6.41%  '
$ws.Range('D35').Value = '''0.0750'
$ws.Range('E35').Value = '  -5.49%  '
$ws.Range('E36').Value = '  -4.77%  '
$ws.Range('D37').Value = '''17.07'
$ws.Range('E37').Value = '  -6.90%  '
$ws.Range('D38').Value = '''1.87'
$ws.Range('E38').Value = '  -7.59%  '
$ws.Range('D39').Value = '''0.103'
$ws.Range('E39').Value = '  -8.28%  '
$ws.Range('E40').Value = '  -4.09%  '
$ws.Range('D41').Value = '''4.07'
$ws.Range('E41').Value = '  -6.21%  '
$ws.Range('D42').Value = '''21.41'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '1.972.10'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '''0.0283'
$ws.Range('E45').Value = '  -5.19%  '
$ws.Range('D46').Value = '''3.04'
$ws.Range('E46').Value = '  -7.82%  '
$ws.Range('D47').Value = '''8.75'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').Value = '''76.73'
$ws.Range('D49').Value = '''96.95'
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('D50').Value = '''68.68'
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('D51').Value = '''0.179'
$ws.Range('E51').Value = '  -6.81%  '
